$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the class moves from column C being '-' to the Tornearia block
$ws.Range("C2").Value = "['MEC-2B-Tornearia', -, -, -]"

# Row 3: class moves from column B to column C
$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "['MEC-2B-Tornearia', -, -, -]"

# Row 4: class removed from column B
$ws.Range("B4").Value = "-"

# Row 6: class moves from column B to column F
$ws.Range("B6").Value = "-"
$ws.Range("F6").Value = "[-, 'MEC-2B-Tornearia', -, -]"

# Row 7: class removed from column B
$ws.Range("B7").Value = "-"

# Row 8: class added to column C
$ws.Range("C8").Value = "[-, 'MEC-2B-Tornearia', -, -]"
